$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.077.30"
$ws.Range("E2").Value = "  -0.52%  "

$ws.Range("D3").Value = "1.784.64"
$ws.Range("E3").Value = "  -2.47%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.97"
$ws.Range("E5").Value = "  +0.21%  "

$ws.Range("E6").Value = "  -1.48%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.66"
$ws.Range("E8").Value = "  +1.58%  "

$ws.Range("E9").Value = "  -2.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0711"
$ws.Range("E10").Value = "  -1.61%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0936"
$ws.Range("E11").Value = "  +0.57%  "

$ws.Range("D12").Value = "2.041.73"
$ws.Range("E12").Value = "  -2.40%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.06"
$ws.Range("E13").Value = "  +1.51%  "

$ws.Range("D14").Value = "1.773.18"
$ws.Range("E14").Value = "  -2.65%  "

$ws.Range("D15").Value = "34.023.08"
$ws.Range("E15").Value = "  -0.74%  "

$ws.Range("E16").Value = "  -3.88%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.15"
$ws.Range("E17").Value = "  -4.55%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.88"
$ws.Range("E18").Value = "  -3.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.85"
$ws.Range("E19").Value = "  -2.97%  "

$ws.Range("D20").Value = "0.0₃0786"
$ws.Range("E20").Value = "  -0.63%  "

$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.73"
$ws.Range("E22").Value = "  -4.40%  "

$ws.Range("E23").Value = "  -4.55%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.09"
$ws.Range("E24").Value = "  -3.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.18"
$ws.Range("E25").Value = "  -0.31%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.32"
$ws.Range("E26").Value = "  -2.27%  "

$ws.Range("E27").Value = "  -3.27%  "

$ws.Range("E28").Value = "  -2.42%  "

$ws.Range("E29").Value = "  +0.18%  "

$ws.Range("E30").Value = "  +0.40%  "

$ws.Range("E31").Value = "  -4.30%  "

$ws.Range("E32").Value = "  -4.25%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.52"
$ws.Range("E33").Value = "  -2.17%  "

$ws.Range("E34").Value = "  -4.91%  "

$ws.Range("D35").Value = "1.390.65"
$ws.Range("E35").Value = "  -3.73%  "

$ws.Range("E36").Value = "  -0.21%  "

$ws.Range("E37").Value = "  -2.22%  "

$ws.Range("E38").Value = "  -1.69%  "

$ws.Range("E39").Value = "  +2.27%  "

$ws.Range("E40").Value = "  -0.21%  "

$ws.Range("E41").Value = "  -5.72%  "

$ws.Range("E42").Value = "  -2.78%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "78.25"
$ws.Range("E43").Value = "  -4.47%  "

$ws.Range("E44").Value = "  +14.02%  "

$ws.Range("E45").Value = "  +2.53%  "

$ws.Range("B46").Value = "Kaspa"
$ws.Range("C46").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0498"
$ws.Range("E46").Value = "  -0.03%  "

$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.47"
$ws.Range("E47").Value = "  +4.38%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "107.93"
$ws.Range("E48").Value = "  +0.84%  "

$ws.Range("E49").Value = "  -4.51%  "

$ws.Range("D50").Value = "1.941.33"
$ws.Range("E50").Value = "  -2.39%  "

$ws.Range("E51").Value = "  +0.13%  "

